$d = $word.ActiveDocument

$d.Content.Find.Execute("详细步骤请参阅FortiToken注册环节", $true, $false, $false, $false, $false, $true, 1, $false, "具体步骤请参阅FortiToken注册环节", 2)

$d.Content.Find.Execute("登录SSL VPN前，您需要安装", $true, $false, $false, $false, $false, $true, 1, $false, "在登录SSL VPN之前，您需要安装", 2)

$d.Content.Find.Execute("打开Google Play商店，搜索", $true, $false, $false, $false, $false, $true, 1, $false, "打开Google Play商店搜索", 2)

$d.Content.Find.Execute("安装“FortiToken Mobile”应用程序后", $true, $false, $false, $false, $false, $true, 1, $false, "安装“FortiToken Mobile”应用后", 2)
